$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 3: "Tv3 Extra " -> "Tv2 С4", price 2500 -> 2100, qty 159 -> 158
$ws.Range("A3").Value = "Tv2 С4"
$ws.Range("B3").Value = 2100
$ws.Range("C3").Value = 158

# Add new row 4: "Ck1 White Diamond", 2200, 177
$ws.Range("A4").Value = "Ck1 White Diamond"
$ws.Range("B4").Value = 2200
$ws.Range("C4").Value = 177

# Update selection to C4 like in the final workbook
$ws.Range("C4").Select() | Out-Null
